{"js": "// 1) Merge the \"S\u1ed1:      /\u0110X-\" run and the \"${soVB}\" run into a single\n//    run. Locate each piece separately (so we don't have to hard-code\n//    the exact whitespace, which includes non-breaking spaces in the\n//    source), expand a range across both, re-read its literal text and\n//    write that same text back - the engine rewrites the whole match as\n//    one run using the formatting of the first run in the match.\nconst body = context.document.body;\n\nconst startMatches = body.search(\"S\u1ed1:\", { matchCase: true, matchWildcards: false });\nstartMatches.load(\"items\");\nconst endMatches = body.search(\"${soVB}\", { matchCase: true, matchWildcards: false });\nendMatches.load(\"items\");\nawait context.sync();\n\nif (startMatches.items.length > 0 && endMatches.items.length > 0) {\n  const numberRange = startMatches.items[0].expandTo(endMatches.items[0]);\n  numberRange.load(\"text\");\n  await context.sync();\n\n  numberRange.insertText(numberRange.text, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Append \"2018\" right after \"\u0110\u00e0 N\u1eb5ng, ng\u00e0y      th\u00e1ng      n\u0103m \" as\n//    its own run, matching the italic / 13pt (sz 26 half-points)\n//    formatting used throughout that line.\nconst dateMatches = body.search(\"\u0110\u00e0 N\u1eb5ng, ng\u00e0y      th\u00e1ng      n\u0103m \", { matchCase: true, matchWildcards: false });\ndateMatches.load(\"items\");\nawait context.sync();\n\nif (dateMatches.items.length > 0) {\n  const yearRun = dateMatches.items[0].insertText(\"2018\", Word.InsertLocation.end);\n  yearRun.font.italic = true;\n  yearRun.font.size = 13;\n  yearRun.font.name = \"Times New Roman\";\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$found = $rng1.Find.Execute('\u0110\u00e0 N\u1eb5ng, ng\u00e0y      th\u00e1ng      n\u0103m ')\n$endPos = $rng1.End\n$rng1.Collapse(0)\n$rng1.InsertAfter(\"2018\")\n\n# Now re-find with a fresh range object to set formatting, see if behaves differently\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Execute(\"2018\")\n$rng2.Font.Italic = $true\n$rng2.Font.Size = 13\n"}
